$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("09/11/2025 00:06:05", "Fundo emergencial ", "Guarde 300", "outros", 300),
    @("09/11/2025 00:06:44", "Saída", "Gastei 550", "outros", 550),
    @("09/11/2025 00:08:00", "Saída", "Gastei 300", "outros", 300),
    @("09/11/2025 00:08:53", "Saída", "Gastei 500", "outros", 500),
    @("11/11/2025 07:16:24", "Entrada", "Recebi 7000", "outros", 7000),
    @("11/11/2025 07:16:55", "Entrada", "Recebi 2000", "outros", 2000)
)

$startRow = 17
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
}
